$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jamal Murray"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Denver Nuggets"
$ws.Range("A3").Value = "CJ McCollum"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "New Orleans Pelicans"
$ws.Range("A4").Value = "Buddy Hield"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "Golden State Warriors"
$ws.Range("A5").Value = "Dillon Brooks"
$ws.Range("B5").Value = "SG,SF,PF"
$ws.Range("C5").Value = "Houston Rockets"
$ws.Range("A6").Value = "Naji Marshall"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Dallas Mavericks"
$ws.Range("A7").Value = "Tobias Harris"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Detroit Pistons"
$ws.Range("A8").Value = "Zach LaVine"
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Sacramento Kings"
$ws.Range("A9").Value = "Kyle Kuzma"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Milwaukee Bucks"
$ws.Range("A10").Value = "Bam Adebayo"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Miami Heat"
$ws.Range("A11").Value = "Donovan Clingan"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Portland Trail Blazers"
$ws.Range("A12").Value = "Kris Dunn"
$ws.Range("B12").Value = "PG,SG"
$ws.Range("C12").Value = "LA Clippers"
$ws.Range("A13").Value = "Khris Middleton"
$ws.Range("B13").Value = "SF"
$ws.Range("C13").Value = "Washington Wizards"
$ws.Range("A14").Value = "John Collins"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Utah Jazz"
$ws.Range("A15").Value = "Shai Gilgeous-Alexander"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Oklahoma City Thunder"
$ws.Range("A16").Value = "Cason Wallace"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Oklahoma City Thunder"
$ws.Range("A17").Value = "Jordan Clarkson"
$ws.Range("B17").Value = "SG,SF"
$ws.Range("C17").Value = "Utah Jazz"
$ws.Range("A18").Value = "Lauri Markkanen"
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Utah Jazz"
$ws.Range("A19").Value = "Jordan Poole"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Washington Wizards"
